$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four "Resolving-Mac" sending-cluster rows (rows 10-13); the
# table now spans A1:T9 instead of A1:T13.
$ws.Rows("10:13").Delete()

# Refresh the remaining TPM-derived metrics (columns G:T) for the eight
# surviving rows (2-9). Columns A:F, K and L are unchanged.
$newValues = @{
    2  = @(0.2899373333333333, 0.869812, 0.5062291280850276, 0.5062291280850276, 18.98824366666667, 56.964731, 0.3642588803316547, 0.3642588803316547, 5.505400733396889, 49.548606600572, 0.184398455387522, 0.1843984553875219)
    3  = @(0.2899373333333333, 0.869812, 0.5062291280850276, 0.5062291280850276, 24.26158266666667, 72.78474800000001, 0.465419398043004, 0.4654193980430039, 7.034338580819556, 63.30904722737601, 0.2356088560651683, 0.2356088560651683)
    4  = @(0.2899373333333333, 0.869812, 0.5062291280850276, 0.5062291280850276, 6.909617666666667, 20.728853, 0.1325498892347874, 0.1325498892347874, 2.003356120626222, 18.030205085636, 0.06710061485509344, 0.06710061485509343)
    5  = @(0.2899373333333333, 0.869812, 0.5062291280850276, 0.5062291280850276, 1.968986333333334, 5.906959000000001, 0.03777183239055392, 0.03777183239055391, 0.5708826468564445, 5.137943821708, 0.01912120177724392, 0.01912120177724391)
    6  = @(0.282802, 0.848406, 0.4937708719149724, 0.4937708719149724, 18.98824366666667, 56.964731, 0.3642588803316547, 0.3642588803316547, 5.369913285420666, 48.329219568786, 0.1798604249441327, 0.1798604249441327)
    7  = @(0.282802, 0.848406, 0.4937708719149724, 0.4937708719149724, 24.26158266666667, 72.78474800000001, 0.465419398043004, 0.4654193980430039, 6.861224101298667, 61.751016911688, 0.2298105419778357, 0.2298105419778356)
    8  = @(0.282802, 0.848406, 0.4937708719149724, 0.4937708719149724, 6.909617666666667, 20.728853, 0.1325498892347874, 0.1325498892347874, 1.954053695368667, 17.586483258318, 0.06544927437969401, 0.06544927437969399)
    9  = @(0.282802, 0.848406, 0.4937708719149724, 0.4937708719149724, 1.968986333333334, 5.906959000000001, 0.03777183239055392, 0.03777183239055391, 0.5568332730393334, 5.011499457354001, 0.01865063061331001, 0.01865063061331)
}

$cols = @("G","H","I","J","M","N","O","P","Q","R","S","T")

foreach ($r in $newValues.Keys) {
    $vals = $newValues[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}
